$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 131.3158
$ws.Range("I5").Value = 131.3158
$ws.Range("K5").Value = 131.3158
$ws.Range("M5").Value = -16.3158
$ws.Range("H6").Value = 3466.6667
$ws.Range("I6").Value = 4000
$ws.Range("K6").Value = 12000
$ws.Range("M6").Value = -11888
$ws.Range("H33").Value = 300.0476
$ws.Range("I33").Value = 137.71428
$ws.Range("K33").Value = 137.71428
$ws.Range("M33").Value = 91.28572
$ws.Range("H70").Value = 7341.2
$ws.Range("I70").Value = 2500
$ws.Range("K70").Value = 7500
$ws.Range("M70").Value = -7230
$ws.Range("H73").Value = 7341.2
$ws.Range("I73").Value = 2500
$ws.Range("K73").Value = 7500
$ws.Range("M73").Value = -6564
$ws.Range("H135").Value = 5409.75
$ws.Range("I135").Value = 5364.75
$ws.Range("K135").Value = 48282.75
$ws.Range("M135").Value = -45747.75
$ws.Range("H137").Value = 15188.863
$ws.Range("J137").Value = 1798.4286
$ws.Range("L137").Value = 5395.2858
$ws.Range("N137").Value = -10495.2858
$ws.Range("H138").Value = 28639.842
$ws.Range("I138").Value = 1781.52
$ws.Range("J138").Value = 80290.46000000001
$ws.Range("K138").Value = 5344.559999999999
$ws.Range("L138").Value = 240871.38
$ws.Range("M138").Value = -204.5599999999995
$ws.Range("N138").Value = -251151.38
$ws.Range("H141").Value = 4966.2
$ws.Range("I141").Value = 5408
$ws.Range("J141").Value = 3199
$ws.Range("K141").Value = 16224
$ws.Range("L141").Value = 9597
$ws.Range("M141").Value = -11044
$ws.Range("N141").Value = -19957

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20464.23
$ws.Range("I32").Value = 20806.666
$ws.Range("K32").Value = 20806.666
$ws.Range("M32").Value = -20519.666
$ws.Range("H61").Value = 6419.6665
$ws.Range("I61").Value = 1087.6
$ws.Range("J61").Value = 19749.834
$ws.Range("K61").Value = 1087.6
$ws.Range("L61").Value = 19749.834
$ws.Range("M61").Value = -875.5999999999999
$ws.Range("N61").Value = -20173.834
$ws.Range("H63").Value = 2132
$ws.Range("J63").Value = 5000
$ws.Range("L63").Value = 5000
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 2132
$ws.Range("J66").Value = 5000
$ws.Range("L66").Value = 25000
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 266674.53
$ws.Range("I74").Value = 429235.94
$ws.Range("K74").Value = 429235.94
$ws.Range("M74").Value = -428361.94
$ws.Range("H77").Value = 266674.53
$ws.Range("I77").Value = 429235.94
$ws.Range("K77").Value = 2146179.7
$ws.Range("M77").Value = -2141811.7
$ws.Range("H122").Value = 2750
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050
$ws.Range("H136").Value = 6419.6665
$ws.Range("I136").Value = 1087.6
$ws.Range("J136").Value = 19749.834
$ws.Range("K136").Value = 3262.8
$ws.Range("L136").Value = 59249.50199999999
$ws.Range("M136").Value = -712.7999999999997
$ws.Range("N136").Value = -64349.50199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1337.5454
$ws.Range("J64").Value = 1185.125
$ws.Range("L64").Value = 1185.125
$ws.Range("N64").Value = -1635.125
$ws.Range("H67").Value = 1337.5454
$ws.Range("J67").Value = 1185.125
$ws.Range("L67").Value = 1185.125
$ws.Range("N67").Value = -2745.125
$ws.Range("H99").Value = 1516.8
$ws.Range("I99").Value = 1309.7142
$ws.Range("K99").Value = 1309.7142
$ws.Range("M99").Value = 188.2858000000001
$ws.Range("H105").Value = 2600.65
$ws.Range("I105").Value = 1946.6154
$ws.Range("K105").Value = 1946.6154
$ws.Range("M105").Value = -199.6153999999999
$ws.Range("H107").Value = 3756.2415
$ws.Range("I107").Value = 4184.5625
$ws.Range("K107").Value = 4184.5625
$ws.Range("M107").Value = -2264.5625
$ws.Range("H134").Value = 2965.1765
$ws.Range("I134").Value = 2426.9092
$ws.Range("K134").Value = 7280.7276
$ws.Range("M134").Value = -4745.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 132.63637
$ws.Range("I7").Value = 73
$ws.Range("J7").Value = 291.66666
$ws.Range("K7").Value = 73
$ws.Range("L7").Value = 291.66666
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = -517.66666
$ws.Range("H31").Value = 5001296
$ws.Range("I31").Value = 5556440
$ws.Range("K31").Value = 5556440
$ws.Range("M31").Value = -5556145
$ws.Range("H34").Value = 5001296
$ws.Range("I34").Value = 5556440
$ws.Range("K34").Value = 5556440
$ws.Range("M34").Value = -5556238
$ws.Range("H94").Value = 3551.2
$ws.Range("J94").Value = 4748.1665
$ws.Range("L94").Value = 4748.1665
$ws.Range("N94").Value = -5650.1665
$ws.Range("H99").Value = 3513.5
$ws.Range("I99").Value = 2199.3333
$ws.Range("K99").Value = 2199.3333
$ws.Range("M99").Value = -701.3332999999998
$ws.Range("H105").Value = 1555.5
$ws.Range("I105").Value = 1492
$ws.Range("K105").Value = 1492
$ws.Range("M105").Value = 255
$ws.Range("H107").Value = 1410.875
$ws.Range("J107").Value = 1637.3334
$ws.Range("L107").Value = 1637.3334
$ws.Range("N107").Value = -5477.3334
$ws.Range("H126").Value = 3513.5
$ws.Range("I126").Value = 2199.3333
$ws.Range("K126").Value = 6597.999899999999
$ws.Range("M126").Value = -4127.999899999999
$ws.Range("H132").Value = 102398.8
$ws.Range("I132").Value = 250752.5
$ws.Range("K132").Value = 752257.5
$ws.Range("M132").Value = -749727.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 7368.421
$ws.Range("H66").Value = 7368.421

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5575
$ws.Range("I122").Value = 5112.5
$ws.Range("K122").Value = 15337.5
$ws.Range("M122").Value = -12887.5
$ws.Range("H126").Value = 2344.875
$ws.Range("I126").Value = 1731.8334
$ws.Range("K126").Value = 5195.5002
$ws.Range("M126").Value = -2725.5002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1672.1818
$ws.Range("I22").Value = 1542
$ws.Range("K22").Value = 1542
$ws.Range("M22").Value = -1247
$ws.Range("H27").Value = 1672.1818
$ws.Range("I27").Value = 1542
$ws.Range("K27").Value = 1542
$ws.Range("M27").Value = -1435
$ws.Range("H40").Value = 2983.6667
$ws.Range("I40").Value = 2979.4
$ws.Range("K40").Value = 2979.4
$ws.Range("M40").Value = -2843.4
$ws.Range("H46").Value = 3686.2942
$ws.Range("I46").Value = 1463.125
$ws.Range("K46").Value = 1463.125
$ws.Range("M46").Value = -1275.125
$ws.Range("H68").Value = 3776.4546
$ws.Range("I68").Value = 3748.2
$ws.Range("K68").Value = 3748.2
$ws.Range("M68").Value = -2999.2
$ws.Range("H71").Value = 3776.4546
$ws.Range("I71").Value = 3748.2
$ws.Range("K71").Value = 18741
$ws.Range("M71").Value = -14997
$ws.Range("H82").Value = 1660.625
$ws.Range("I82").Value = 1588.4546
$ws.Range("K82").Value = 1588.4546
$ws.Range("M82").Value = -1227.4546
$ws.Range("H85").Value = 1660.625
$ws.Range("I85").Value = 1588.4546
$ws.Range("K85").Value = 1588.4546
$ws.Range("M85").Value = -340.4546
$ws.Range("H132").Value = 3446.3684
$ws.Range("I132").Value = 2989.4482
$ws.Range("J132").Value = 4918.6665
$ws.Range("K132").Value = 8968.3446
$ws.Range("L132").Value = 14755.9995
$ws.Range("M132").Value = -6438.3446
$ws.Range("N132").Value = -19815.9995
$ws.Range("H136").Value = 4144.8237
$ws.Range("I136").Value = 3757.5
$ws.Range("J136").Value = 4489.1113
$ws.Range("K136").Value = 11272.5
$ws.Range("L136").Value = 13467.3339
$ws.Range("M136").Value = -8722.5
$ws.Range("N136").Value = -18567.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 728.63635
$ws.Range("I100").Value = 829.2857
$ws.Range("K100").Value = 1658.5714
$ws.Range("M100").Value = -1117.5714
$ws.Range("H132").Value = 2211.7878
$ws.Range("I132").Value = 1633.2727
$ws.Range("K132").Value = 4899.8181
$ws.Range("M132").Value = -2369.8181
